# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-12-03 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-12-04 Thursday", 2)

$t = $d.Tables.Item(1)

# Row 1: first cell's problem is dropped, the remaining four shift left
# (keeping the unchanged "24/2=12,0" in the new first slot) and a new
# problem is appended as the new last cell.
$t.Cell(1, 1).Range.Text = "24÷2=12, 0"
$t.Cell(1, 2).Range.Text = "52÷6=8, 4"
$t.Cell(1, 3).Range.Text = "49÷8=6, 1"
$t.Cell(1, 4).Range.Text = "14÷7=2, 0"
$t.Cell(1, 5).Range.Text = "18÷3=6, 0"

# Row 5: straight text swaps.
$t.Cell(5, 1).Range.Text = "25÷7=3, 4"
$t.Cell(5, 2).Range.Text = "51÷5=10, 1"
$t.Cell(5, 3).Range.Text = "16÷4=4, 0"
$t.Cell(5, 4).Range.Text = "91÷4=22, 3"
$t.Cell(5, 5).Range.Text = "72÷2=36, 0"

# Row 9: straight text swaps.
$t.Cell(9, 1).Range.Text = "94÷3=31, 1"
$t.Cell(9, 2).Range.Text = "64÷5=12, 4"
$t.Cell(9, 3).Range.Text = "80÷3=26, 2"
$t.Cell(9, 4).Range.Text = "41÷8=5, 1"
$t.Cell(9, 5).Range.Text = "82÷7=11, 5"

# Row 13: straight text swaps.
$t.Cell(13, 1).Range.Text = "78÷8=9, 6"
$t.Cell(13, 2).Range.Text = "77÷7=11, 0"
$t.Cell(13, 3).Range.Text = "34÷8=4, 2"
$t.Cell(13, 4).Range.Text = "18÷2=9, 0"
$t.Cell(13, 5).Range.Text = "80÷2=40, 0"

# Row 17: straight text swaps.
$t.Cell(17, 1).Range.Text = "15÷9=1, 6"
$t.Cell(17, 2).Range.Text = "80÷8=10, 0"
$t.Cell(17, 3).Range.Text = "10÷4=2, 2"
$t.Cell(17, 4).Range.Text = "68÷2=34, 0"
$t.Cell(17, 5).Range.Text = "78÷4=19, 2"
